$d = $word.ActiveDocument

function ReplaceText($findText, $replaceText) {
    $range = $d.Content
    $ok = $range.Find.Execute($findText, $true, $false, $false, $false, $false, `
                               $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $findText"
    }
}

# 1. "Python Programming" skill line: merge " I " + "Jupyter" + " Notebook" into one run
ReplaceText " I Jupyter Notebook" " I Jupyter Notebook"
# ... and merge " | " + "Sklearn" + " | " + "Xgboost" + " | Bs4" into one run
ReplaceText " | Sklearn | Xgboost | Bs4" " | Sklearn | Xgboost | Bs4"

# 2. Databases line: append " | AWS (RDS, S3)" after "MySQL | MongoDB | ETL"
ReplaceText "MySQL | MongoDB | ETL" "MySQL | MongoDB | ETL | AWS (RDS, S3)"

# 3. Web Technologies line: merge the Geomapping runs into one
ReplaceText "HTML | CSS | Bootstrap | Dashboarding | JavaScript Charting | D3.js | Geomapping with Leaflet.js" `
            "HTML | CSS | Bootstrap | Dashboarding | JavaScript Charting | D3.js | Geomapping with Leaflet.js"

# 4. Covid-Prediction project skill line: merge Sklearn/Xgboost runs
ReplaceText "Python | Pandas | Seaborn | Sklearn | Matplotlib | Xgboost | Bs4" `
            "Python | Pandas | Seaborn | Sklearn | Matplotlib | Xgboost | Bs4"

# 5. Second project skill line: merge Sklearn run
ReplaceText "Python | Pandas | Sklearn | Matplotlib" "Python | Pandas | Sklearn | Matplotlib"

# 6. NBA-Cluster project title: merge "NBA-" + "Cluster" + " " into one run
ReplaceText "NBA-Cluster " "NBA-Cluster "

# 7. NBA clustering skill line: merge " | NumPy | Matplotlib | " + "Sklearn" + " " into one run
ReplaceText " | NumPy | Matplotlib | Sklearn " " | NumPy | Matplotlib | Sklearn "

# 8. Assistant Project Manager bullet: merge "Managed subcontractors " + "which" + rest into one run
ReplaceText "Managed subcontractors which ultimately expedited the construction of different custom houses" `
            "Managed subcontractors which ultimately expedited the construction of different custom houses"
